$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "timestamp" column (O) for every data row (2 through 530)
# from the old crawl time "2023-01-15 06:50:55" to the new crawl time
# "2023-01-15 12:56:47".
$ws.Range("O2:O530").Value = "2023-01-15 12:56:47"
